$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Update "Pagos" (F) and "Inscrições homologadas" (H) columns
# Row 2: 10 -> 12
$ws.Range("F2").Value = 12
$ws.Range("H2").Value = 12

# Row 3: 9 -> 10
$ws.Range("F3").Value = 10
$ws.Range("H3").Value = 10

# Row 7: 13 -> 14
$ws.Range("F7").Value = 14
$ws.Range("H7").Value = 14

# Row 10: 10 -> 11
$ws.Range("F10").Value = 11
$ws.Range("H10").Value = 11

# Row 11: 9 -> 10
$ws.Range("F11").Value = 10
$ws.Range("H11").Value = 10

# Row 12: 8 -> 9
$ws.Range("F12").Value = 9
$ws.Range("H12").Value = 9

# Row 16: 82 -> 84
$ws.Range("F16").Value = 84
$ws.Range("H16").Value = 84

$wb.Save()
